$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression
$ws.Range("B2").Value = 0.6994913607464788
$ws.Range("C2").Value = 0.6994913607464788
$ws.Range("D2").Value = 0.6994913607464788

# Row 3: RandomForestRegressor
$ws.Range("B3").Value = 0.9790096649213884
$ws.Range("C3").Value = 0.9795082789696985
$ws.Range("D3").Value = 0.9793517347724484

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9812305404534789
$ws.Range("C4").Value = 0.9813911928860608
$ws.Range("D4").Value = 0.9811164696555958

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8343912838367032
$ws.Range("C5").Value = 0.7930020961697489
$ws.Range("D5").Value = 0.8023159764513645
